$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) for new columns I and J, matching style of existing headers (bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-26 for columns I (I0) and J (IF)
$data = @{
    2  = @(7, 7)
    3  = @(9, 9)
    4  = @(7, 8)
    5  = @(8, 8)
    6  = @(6, 6)
    7  = @(6, 7)
    8  = @(6, 7)
    9  = @(6, 7)
    10 = @(6, 7)
    11 = @(7, 8)
    12 = @(1, 5)
    13 = @(1, 8)
    14 = @(1, 6)
    15 = @(1, 1)
    16 = @(1, 5)
    17 = @(1, 3)
    18 = @(1, 4)
    19 = @(1, 6)
    20 = @(1, 4)
    21 = @(1, 5)
    22 = @(8, 8)
    23 = @(1, 6)
    24 = @(1, 5)
    25 = @(5, 7)
    26 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
